$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column G re-shuffle (the "sub-optimal ACL generation") ---
# Row 2's "Reject" (G) cell ends up re-pointed to a newly-inserted string
# instead of being cleared like the rest - that's the known bug being
# committed here.
$ws.Range("G2").Value = "rogue"

# Rows 3-6, 8 and 9 lose their "ALL" Reject entry entirely.
$ws.Range("G3").Value = ""
$ws.Range("G4").Value = ""
$ws.Range("G5").Value = ""
$ws.Range("G6").Value = ""
$ws.Range("G8").Value = ""
$ws.Range("G9").Value = ""

# Row 6's Status flips from disable to enable.
$ws.Range("H6").Value = "enable"

# Row 8: Accept becomes ALL (was PT), Reject cleared (above), Status flips to enable.
$ws.Range("F8").Value = "ALL"
$ws.Range("H8").Value = "enable"

# --- Sheet view: move the active selection ---
[void]$ws.Range("F17").Select()

# --- Column widths ---
$ws.Columns.Item(7).ColumnWidth = 36.6
$ws.Columns.Item(8).ColumnWidth = 17
